$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7228397727012634
$ws.Range("B1").Value = 1.097208738327026
$ws.Range("C1").Value = 2.403980731964111
$ws.Range("D1").Value = 3.616001844406128
$ws.Range("E1").Value = 1.719799757003784
